$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Summary header updates ---
$ws.Range("E11").Value = 1065654      # VALOR MORA total
$ws.Range("C13").Value = 2            # Cant. Trabajadores
$ws.Range("F13").Value = 17           # Cant. Periodos

# --- Detail rows 16-38 (two workers: Carlos Andres Perez Martinez, Marta Isabel Ledezma Pabuena) ---
$data = @(
  @(16,"CC","1047442827","CARLOS ANDRES PEREZ MARTINEZ","2303",46400,1160000),
  @(17,"CC","1047442827","CARLOS ANDRES PEREZ MARTINEZ","2304",46400,1160000),
  @(18,"CC","33207649","MARTA ISABEL LEDEZMA PABUENA","2304",46400,1160000),
  @(19,"CC","1047442827","CARLOS ANDRES PEREZ MARTINEZ","2305",46400,1160000),
  @(20,"CC","33207649","MARTA ISABEL LEDEZMA PABUENA","2305",46400,1160000),
  @(21,"CC","1047442827","CARLOS ANDRES PEREZ MARTINEZ","2306",46400,1160000),
  @(22,"CC","33207649","MARTA ISABEL LEDEZMA PABUENA","2306",46400,1160000),
  @(23,"CC","1047442827","CARLOS ANDRES PEREZ MARTINEZ","2307",46400,1160000),
  @(24,"CC","33207649","MARTA ISABEL LEDEZMA PABUENA","2307",46400,1160000),
  @(25,"CC","1047442827","CARLOS ANDRES PEREZ MARTINEZ","2308",46400,1160000),
  @(26,"CC","33207649","MARTA ISABEL LEDEZMA PABUENA","2308",46400,1160000),
  @(27,"CC","1047442827","CARLOS ANDRES PEREZ MARTINEZ","2309",46400,1160000),
  @(28,"CC","33207649","MARTA ISABEL LEDEZMA PABUENA","2309",46400,1160000),
  @(29,"CC","1047442827","CARLOS ANDRES PEREZ MARTINEZ","2310",46400,1160000),
  @(30,"CC","1047442827","CARLOS ANDRES PEREZ MARTINEZ","2311",46400,1160000),
  @(31,"CC","1047442827","CARLOS ANDRES PEREZ MARTINEZ","2312",46400,1160000),
  @(32,"CC","1047442827","CARLOS ANDRES PEREZ MARTINEZ","2401",46400,1160000),
  @(33,"CC","1047442827","CARLOS ANDRES PEREZ MARTINEZ","2402",46400,1160000),
  @(34,"CC","1047442827","CARLOS ANDRES PEREZ MARTINEZ","2403",46400,1160000),
  @(35,"CC","1047442827","CARLOS ANDRES PEREZ MARTINEZ","2404",46400,1160000),
  @(36,"CC","1047442827","CARLOS ANDRES PEREZ MARTINEZ","2405",46400,1160000),
  @(37,"CC","1047442827","CARLOS ANDRES PEREZ MARTINEZ","2406",46400,1160000),
  @(38,"CC","1047442827","CARLOS ANDRES PEREZ MARTINEZ","2407",44854,1160000)
)

foreach ($row in $data) {
  $r = $row[0]
  $ws.Range("B$r").Value = $row[1]
  $ws.Range("C$r").Value = $row[2]
  $ws.Range("D$r").Value = $row[3]
  $ws.Range("E$r").Value = $row[4]
  $ws.Range("F$r").Value = $row[5]
  $ws.Range("G$r").Value = $row[6]
}

# Row 38 becomes the new last row of the table -- carry over the closing
# bottom-border formatting that used to live on row 41 before we remove the
# now-stale rows 39:41.
$ws.Range("B41:J41").Copy()
$ws.Range("B38:J38").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false

$ws.Rows("39:41").Delete()

# Column D ("Nombre Trabajador") no longer needs to fit the longer names that
# were removed -- re-fit it against the remaining (shorter) content.
$ws.Columns("D").ColumnWidth = 31.6
